$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff cyclically rotates the data (Fecha/Volumen/Precios/Unidad/Origen/etc.)
# held in rows 3, 4 and 5 of columns D and M:T:
#   new row3 = old row5
#   new row4 = old row3
#   new row5 = old row4
# Columns A:C, E:L stay identical across these rows, so only D and M:T need updating.

# Capture the "before" values first, since we are about to overwrite them.
# (use Value2 - Value getter is unreliable in this runtime)
$old3_D = $ws.Range("D3").Value2
$old3_M = $ws.Range("M3").Value2
$old3_N = $ws.Range("N3").Value2
$old3_O = $ws.Range("O3").Value2
$old3_P = $ws.Range("P3").Value2
$old3_Q = $ws.Range("Q3").Value2
$old3_R = $ws.Range("R3").Value2
$old3_S = $ws.Range("S3").Value2
$old3_T = $ws.Range("T3").Value2

$old4_D = $ws.Range("D4").Value2
$old4_M = $ws.Range("M4").Value2
$old4_N = $ws.Range("N4").Value2
$old4_O = $ws.Range("O4").Value2
$old4_P = $ws.Range("P4").Value2
$old4_Q = $ws.Range("Q4").Value2
$old4_R = $ws.Range("R4").Value2
$old4_S = $ws.Range("S4").Value2
$old4_T = $ws.Range("T4").Value2

$old5_D = $ws.Range("D5").Value2
$old5_M = $ws.Range("M5").Value2
$old5_N = $ws.Range("N5").Value2
$old5_O = $ws.Range("O5").Value2
$old5_P = $ws.Range("P5").Value2
$old5_Q = $ws.Range("Q5").Value2
$old5_R = $ws.Range("R5").Value2
$old5_S = $ws.Range("S5").Value2
$old5_T = $ws.Range("T5").Value2

# Row 3 <- old row 5
$ws.Range("D3").Value2 = $old5_D
$ws.Range("M3").Value2 = $old5_M
$ws.Range("N3").Value2 = $old5_N
$ws.Range("O3").Value2 = $old5_O
$ws.Range("P3").Value2 = $old5_P
$ws.Range("Q3").Value2 = $old5_Q
$ws.Range("R3").Value2 = $old5_R
$ws.Range("S3").Value2 = $old5_S
$ws.Range("T3").Value2 = $old5_T

# Row 4 <- old row 3
$ws.Range("D4").Value2 = $old3_D
$ws.Range("M4").Value2 = $old3_M
$ws.Range("N4").Value2 = $old3_N
$ws.Range("O4").Value2 = $old3_O
$ws.Range("P4").Value2 = $old3_P
$ws.Range("Q4").Value2 = $old3_Q
$ws.Range("R4").Value2 = $old3_R
$ws.Range("S4").Value2 = $old3_S
$ws.Range("T4").Value2 = $old3_T

# Row 5 <- old row 4
$ws.Range("D5").Value2 = $old4_D
$ws.Range("M5").Value2 = $old4_M
$ws.Range("N5").Value2 = $old4_N
$ws.Range("O5").Value2 = $old4_O
$ws.Range("P5").Value2 = $old4_P
$ws.Range("Q5").Value2 = $old4_Q
$ws.Range("R5").Value2 = $old4_R
$ws.Range("S5").Value2 = $old4_S
$ws.Range("T5").Value2 = $old4_T
